$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in B6 and C6 but keep their existing formatting/style
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# Update the selection to D16
$ws.Range("D16").Select()
